$d = $word.ActiveDocument

# 1. Title heading + the bolded title reprise later in the document (identical text, both occurrences change identically)
$d.Content.Find.Execute(
    "Play Drift King Free: A Unique Game with High-speed Thrill", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Drift King Free: Exciting Racing Slot Game", 2)

# 2. "What we like" bullets
$d.Content.Find.Execute(
    "High-tech racing visuals", $true, $false, $false, $false, $false,
    $true, 1, $false, "High-tech visuals and immersive racing theme", 2)

$d.Content.Find.Execute(
    "Above-average RTP", $true, $false, $false, $false, $false,
    $true, 1, $false, "Above-average RTP for increased profitability", 2)

$d.Content.Find.Execute(
    "Multiple bonus features", $true, $false, $false, $false, $false,
    $true, 1, $false, "Multiple bonus features for higher winnings", 2)

$d.Content.Find.Execute(
    "Challenging gameplay", $true, $false, $false, $false, $false,
    $true, 1, $false, "Challenging gameplay with exciting tension", 2)

# 3. "What we don't like" bullets - swap + reword, done via direct paragraph
#    text assignment (paragraphs 46 and 47) since the two bullets collide
#    under naive global Find/Replace.
$p46 = $d.Paragraphs(46).Range
$p46.Find.Execute("Disruptive music", $true, $false, $false, $false, $false,
                   $true, 1, $false, "Potentially confusing color palette", 2)

$p47 = $d.Paragraphs(47).Range
$p47.Find.Execute("Potentially confusing color palette", $true, $false, $false, $false, $false,
                   $true, 1, $false, "Constant music may be disruptive", 2)

# 4. Meta description (italic) paragraph near the end of the document
$d.Content.Find.Execute(
    "Read our review of Drift King, a unique racing game with multiple bonuses, and play for free. Enjoy challenging gameplay and immerse in high-speed visuals.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Read our review of Drift King, a thrilling racing-themed slot game. Play for free and enjoy challenging gameplay.",
    2)
